$d = $word.ActiveDocument

# Locate the final run of the paragraph: ". The board is round and is powered on a battery."
# and remove it entirely so we can rebuild it as several distinct runs, matching
# the structure produced by the original authoring/editing session.
$old = ". The board is round and is powered on a battery."
$rng = $d.Content
$found = $rng.Find.Execute($old)
if (-not $found) {
    throw "Could not find the target sentence to edit."
}
$rng.Delete()

# Helper: insert a new run of text at the very end of the document content.
function Append-Run([string]$text) {
    $end = $d.Content.End
    $ins = $d.Range($end, $end)
    $ins.InsertAfter($text)
}

Append-Run " in an up/down/left/right/center configuration"
Append-Run ". "
Append-Run "The user can also create custom emojis"
Append-Run " and play games"
Append-Run " using these buttons. "
Append-Run "The board is round and is powered on a battery."
